$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to text format so numeric-looking strings
# (e.g. "43.09", "68.212.28") are preserved as text, not converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "68.212.28"
$ws.Range("E2").Value = "  +1.94%  "
$ws.Range("D3").Value = "3.910.83"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").Value = "483.51"
$ws.Range("E5").Value = "  +3.13%  "
$ws.Range("D6").Value = "146.70"
$ws.Range("E6").Value = "  +1.89%  "
$ws.Range("D7").Value = "0.620"
$ws.Range("E7").Value = "  +1.70%  "
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").Value = "0.736"
$ws.Range("E9").Value = "  +3.79%  "
$ws.Range("D10").Value = "0.169"
$ws.Range("E10").Value = "  +4.19%  "
$ws.Range("D11").Value = "0.0000345"
$ws.Range("E11").Value = "  +2.38%  "
$ws.Range("D12").Value = "43.09"
$ws.Range("E12").Value = "  +2.91%  "
$ws.Range("D13").Value = "10.74"
$ws.Range("E13").Value = "  +6.45%  "
$ws.Range("D14").Value = "4.538.13"
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("D15").Value = "3.912.66"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("D16").Value = "14.20"
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "20.20"
$ws.Range("E18").Value = "  +3.69%  "
$ws.Range("E19").Value = "  +2.81%  "
$ws.Range("D20").Value = "68.302.35"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").Value = "430.53"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("D22").Value = "3.54"
$ws.Range("E22").Value = "  +9.26%  "
$ws.Range("D23").Value = "15.03"
$ws.Range("E23").Value = "  +6.68%  "
$ws.Range("D24").Value = "89.33"
$ws.Range("E24").Value = "  +4.08%  "
$ws.Range("D25").Value = "11.67"
$ws.Range("E25").Value = "  +17.66%  "
$ws.Range("D26").Value = "3.72"
$ws.Range("E26").Value = "  +4.18%  "
$ws.Range("D27").Value = "11.10"
$ws.Range("E27").Value = "  +11.13%  "
$ws.Range("D28").Value = "37.59"
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").Value = "5.68"
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("D30").Value = "717.71"
$ws.Range("E30").Value = "  -1.17%  "
$ws.Range("D31").Value = "13.73"
$ws.Range("E31").Value = "  +5.29%  "
$ws.Range("E32").Value = "  +4.74%  "
$ws.Range("E33").Value = "  +5.11%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "41.92"
$ws.Range("E34").Value = "  +1.51%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "6.10"
$ws.Range("E35").Value = "  +15.59%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0878"
$ws.Range("E36").Value = "  +4.23%  "
$ws.Range("D37").Value = "60.77"
$ws.Range("E37").Value = "  +4.71%  "
$ws.Range("D38").Value = "0.401"
$ws.Range("E38").Value = "  +20.74%  "
$ws.Range("D39").Value = "3.04"
$ws.Range("E39").Value = "  +12.35%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.146"
$ws.Range("E40").Value = "  -2.03%  "
$ws.Range("B41").Value = "Dai"
$ws.Range("C41").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").Value = "0.0493"
$ws.Range("E42").Value = "  +6.77%  "
$ws.Range("D43").Value = "3.12"
$ws.Range("E43").Value = "  +4.93%  "
$ws.Range("E44").Value = "  +3.28%  "
$ws.Range("E45").Value = "  +3.28%  "
$ws.Range("E46").Value = "  +5.97%  "
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("E48").Value = "  +3.20%  "
$ws.Range("D49").Value = "2.13"
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("D50").Value = "145.16"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").Value = "0.0₆0333"
$ws.Range("E51").Value = "  +31.26%  "

# Restore default (Normal) style for the Price/Volume columns so no stray
# text-format styling is left behind on cells.
$ws.Range("D2:E51").Style = "Normal"

